$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "variation array" values in column A (bin edges) and the
# --- "relative frequencies" values in column B that changed alongside them.
$ws.Range("A1").Value = 73.7
$ws.Range("A2").Value = 82.26
$ws.Range("A3").Value = 90.82
$ws.Range("A4").Value = 99.38
$ws.Range("A5").Value = 107.94
$ws.Range("B5").Value = 0.29
$ws.Range("A6").Value = 116.5
$ws.Range("B6").Value = 0.12
$ws.Range("A7").Value = 125.06
$ws.Range("A8").Value = 133.62

# --- Update the histogram chart to match the recalculated data.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$cg = $chart.ChartGroups(1)

# varyColors: true -> false
$cg.VaryByCategories = $false

# gapWidth: add (not present before) with value 2
$cg.GapWidth = 2
